$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.256.77"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.344.31"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'303.11"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "'95.58"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.497"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'34.30"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").Value = "'0.0786"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'18.63"
$ws.Range("E12").Value = "  -3.49%  "
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "2.706.08"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "2.345.82"
$ws.Range("E16").Value = "  +1.51%  "
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "43.182.80"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "'12.27"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'6.22"
$ws.Range("E20").Value = "  +3.52%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "'68.07"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").Value = "'236.19"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("E24").Value = "  -0.43%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").Value = "'2.43"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'24.64"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "'2.35"
$ws.Range("E28").Value = "  +7.11%  "
$ws.Range("D29").Value = "'9.22"
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("D30").Value = "'31.54"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").Value = "'0.0733"
$ws.Range("E33").Value = "  +5.37%  "
$ws.Range("D34").Value = "'17.39"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("E35").Value = "  +5.14%  "
$ws.Range("D36").Value = "'4.37"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'22.53"
$ws.Range("E39").Value = "  +17.57%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "'2.77"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'117.47"
$ws.Range("E42").Value = "  -28.42%  "
$ws.Range("D43").Value = "1.939.89"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").Value = "'10.02"
$ws.Range("E45").Value = "  -4.64%  "
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").Value = "'2.73"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").Value = "2.571.08"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("D51").Value = "'72.25"
$ws.Range("E51").Value = "  +0.12%  "
